# Auto-generated edit script applying odds updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 7

# Row 17
$ws.Range("G17").Value = 1.9
$ws.Range("H17").Value = 3.3
$ws.Range("J17").Value = 2.63
$ws.Range("K17").Value = 2
$ws.Range("S17").Value = 1.53
$ws.Range("T17").Value = 2.38
$ws.Range("U17").Value = 2.1
$ws.Range("V17").Value = 1.67
$ws.Range("X17").Value = 8
$ws.Range("Y17").Value = 9
$ws.Range("AB17").Value = 34
$ws.Range("AC17").Value = 7.5
$ws.Range("AE17").Value = 19
$ws.Range("AF17").Value = 67
$ws.Range("AI17").Value = 19
$ws.Range("AK17").Value = 41
$ws.Range("AS17").Value = 201
$ws.Range("AT17").Value = 2.38
$ws.Range("AU17").Value = 9
$ws.Range("AV17").Value = 67
$ws.Range("BC17").Value = 126
$ws.Range("BD17").Value = 351

# Row 18
$ws.Range("K18").Value = 1.8
$ws.Range("U18").Value = 2.5
$ws.Range("V18").Value = 1.5
$ws.Range("AB18").Value = 51

# Row 19
$ws.Range("J19").Value = 2.5
$ws.Range("K19").Value = 1.91
$ws.Range("L19").Value = 6.5
$ws.Range("S19").Value = 1.62
$ws.Range("T19").Value = 2.2
$ws.Range("U19").Value = 2.63
$ws.Range("V19").Value = 1.44
$ws.Range("W19").Value = 4.75
$ws.Range("Y19").Value = 10
$ws.Range("AC19").Value = 6
$ws.Range("AE19").Value = 26
$ws.Range("AH19").Value = 9.5
$ws.Range("AO19").Value = 10
$ws.Range("AP19").Value = 29
$ws.Range("AS19").Value = 301
$ws.Range("AT19").Value = 2.2
$ws.Range("AU19").Value = 11
$ws.Range("AV19").Value = 101
$ws.Range("BC19").Value = 251

# Row 23
$ws.Range("G23").Value = 3.3
$ws.Range("H23").Value = 2.88
$ws.Range("I23").Value = 2.4
$ws.Range("J23").Value = 4
$ws.Range("L23").Value = 3.25
$ws.Range("M23").Value = 1.13
$ws.Range("N23").Value = 6
$ws.Range("U23").Value = 2.2
$ws.Range("V23").Value = 1.62
$ws.Range("X23").Value = 15
$ws.Range("Y23").Value = 13
$ws.Range("AA23").Value = 34
$ws.Range("AE23").Value = 19
$ws.Range("AH23").Value = 6
$ws.Range("AI23").Value = 10
$ws.Range("AK23").Value = 23
$ws.Range("AL23").Value = 23
$ws.Range("AN23").Value = 5
$ws.Range("AO23").Value = 21
$ws.Range("AS23").Value = 351
$ws.Range("AV23").Value = 81

# Row 29
$ws.Range("G29").Value = 1.55
$ws.Range("H29").Value = 3.6
$ws.Range("I29").Value = 6.25
$ws.Range("Q29").Value = 2.25
$ws.Range("R29").Value = 1.62
$ws.Range("AB29").Value = 34
$ws.Range("AD29").Value = 7.5

# Row 45
$ws.Range("M45").Value = 1.03
$ws.Range("N45").Value = 15
$ws.Range("O45").Value = 1.2
$ws.Range("P45").Value = 4.33

# Row 46
$ws.Range("G46").Value = 1.47
$ws.Range("H46").Value = 4.15
$ws.Range("I46").Value = 5.9
$ws.Range("J46").Value = 1.98
$ws.Range("K46").Value = 2.27
$ws.Range("L46").Value = 5.7
$ws.Range("U46").Value = 1.88
$ws.Range("V46").Value = 1.72
$ws.Range("W46").Value = 6.7
$ws.Range("X46").Value = 6.7
$ws.Range("Z46").Value = 9.75
$ws.Range("AB46").Value = 28
$ws.Range("AC46").Value = 11.25
$ws.Range("AD46").Value = 8.25
$ws.Range("AE46").Value = 19
$ws.Range("AF46").Value = 90
$ws.Range("AG46").Value = 800
$ws.Range("AH46").Value = 15.5
$ws.Range("AI46").Value = 37
$ws.Range("AJ46").Value = 19
$ws.Range("AK46").Value = 120
$ws.Range("AL46").Value = 65
$ws.Range("AM46").Value = 65
$ws.Range("AN46").Value = 3.25
$ws.Range("AO46").Value = 6.8
$ws.Range("AP46").Value = 17.5
$ws.Range("AQ46").Value = 20
$ws.Range("AR46").Value = 50
$ws.Range("AU46").Value = 8.25
$ws.Range("AV46").Value = 80
$ws.Range("AY46").Value = 7.2
$ws.Range("AZ46").Value = 35
$ws.Range("BA46").Value = 37
$ws.Range("BB46").Value = 250

# Row 47
$ws.Range("G47").Value = 1.37
$ws.Range("H47").Value = 4.5
$ws.Range("I47").Value = 7.3
$ws.Range("J47").Value = 1.83
$ws.Range("K47").Value = 2.4
$ws.Range("L47").Value = 6.3
$ws.Range("T47").Value = 3.42
$ws.Range("U47").Value = 1.8
$ws.Range("V47").Value = 1.8
$ws.Range("X47").Value = 6.9
$ws.Range("Z47").Value = 9
$ws.Range("AB47").Value = 24
$ws.Range("AD47").Value = 9.25
$ws.Range("AE47").Value = 18.5
$ws.Range("AF47").Value = 75
$ws.Range("AM47").Value = 65
$ws.Range("AO47").Value = 6.1
$ws.Range("AQ47").Value = 16.5
$ws.Range("AT47").Value = 3.15
$ws.Range("AU47").Value = 7.9
$ws.Range("AZ47").Value = 40

# Row 60
$ws.Range("G60").Value = 2.4
$ws.Range("I60").Value = 3
$ws.Range("J60").Value = 3.25
$ws.Range("K60").Value = 1.95
$ws.Range("L60").Value = 3.75
$ws.Range("M60").Value = 1.08
$ws.Range("N60").Value = 8
$ws.Range("O60").Value = 1.44
$ws.Range("P60").Value = 2.63
$ws.Range("S60").Value = 1.53
$ws.Range("T60").Value = 2.38
$ws.Range("X60").Value = 11
$ws.Range("Z60").Value = 23
$ws.Range("AA60").Value = 23
$ws.Range("AB60").Value = 41
$ws.Range("AF60").Value = 67
$ws.Range("AH60").Value = 8
$ws.Range("AI60").Value = 13
$ws.Range("AJ60").Value = 12
$ws.Range("AO60").Value = 15
$ws.Range("AP60").Value = 29
$ws.Range("AQ60").Value = 51
$ws.Range("AS60").Value = 251
$ws.Range("AT60").Value = 2.38
$ws.Range("AY60").Value = 4.75
$ws.Range("BA60").Value = 29
$ws.Range("BB60").Value = 51

# Row 61
$ws.Range("G61").Value = 4.75
$ws.Range("I61").Value = 1.75
$ws.Range("S61").Value = 1.5
$ws.Range("T61").Value = 2.5
$ws.Range("U61").Value = 2.1
$ws.Range("V61").Value = 1.67
$ws.Range("X61").Value = 23
$ws.Range("Y61").Value = 17
$ws.Range("AB61").Value = 51
$ws.Range("AE61").Value = 19
$ws.Range("AI61").Value = 7.5
$ws.Range("AK61").Value = 13
$ws.Range("AN61").Value = 6.5
$ws.Range("AP61").Value = 41
$ws.Range("AS61").Value = 351
$ws.Range("AT61").Value = 2.5
$ws.Range("AY61").Value = 3.6

# Row 87
$ws.Range("G87").Value = 3.9
$ws.Range("H87").Value = 3.3
$ws.Range("I87").Value = 1.9
$ws.Range("J87").Value = 3.9
$ws.Range("K87").Value = 2.07
$ws.Range("L87").Value = 2.62
$ws.Range("M87").Value = 1.03
$ws.Range("N87").Value = 11.9
$ws.Range("Q87").Value = 1.93
$ws.Range("R87").Value = 1.78
$ws.Range("U87").Value = 1.75
$ws.Range("V87").Value = 1.85
$ws.Range("W87").Value = 10.75
$ws.Range("X87").Value = 21
$ws.Range("Y87").Value = 13
$ws.Range("Z87").Value = 60
$ws.Range("AA87").Value = 37
$ws.Range("AB87").Value = 45
$ws.Range("AD87").Value = 6.4
$ws.Range("AE87").Value = 14.5
$ws.Range("AF87").Value = 70
$ws.Range("AG87").Value = 600
$ws.Range("AH87").Value = 6.9
$ws.Range("AI87").Value = 9
$ws.Range("AJ87").Value = 8.25
$ws.Range("AK87").Value = 16.5
$ws.Range("AL87").Value = 15.5
$ws.Range("AM87").Value = 27
$ws.Range("AN87").Value = 5.4
$ws.Range("AO87").Value = 19
$ws.Range("AP87").Value = 24
$ws.Range("AQ87").Value = 90
$ws.Range("AR87").Value = 120
$ws.Range("AS87").Value = 300
$ws.Range("AY87").Value = 3.9
$ws.Range("AZ87").Value = 10.25

# Row 88
$ws.Range("G88").Value = 3.25
$ws.Range("H88").Value = 2.92
$ws.Range("J88").Value = 4.05
$ws.Range("K88").Value = 1.87
$ws.Range("L88").Value = 2.95
$ws.Range("N88").Value = 6.95
$ws.Range("T88").Value = 2.2
$ws.Range("U88").Value = 1.83
$ws.Range("V88").Value = 1.78
$ws.Range("W88").Value = 8
$ws.Range("X88").Value = 16
$ws.Range("Y88").Value = 11.75
$ws.Range("AA88").Value = 35
$ws.Range("AB88").Value = 45
$ws.Range("AC88").Value = 7.3
$ws.Range("AD88").Value = 5.7
$ws.Range("AG88").Value = 700
$ws.Range("AH88").Value = 6.9
$ws.Range("AI88").Value = 10.75
$ws.Range("AJ88").Value = 9
$ws.Range("AK88").Value = 24
$ws.Range("AL88").Value = 19.5
$ws.Range("AM88").Value = 32
$ws.Range("AN88").Value = 5
$ws.Range("AO88").Value = 20
$ws.Range("AP88").Value = 32
$ws.Range("AQ88").Value = 120
$ws.Range("AR88").Value = 200
$ws.Range("AS88").Value = 500
$ws.Range("AT88").Value = 2.18
$ws.Range("AZ88").Value = 12.5
$ws.Range("BA88").Value = 23
$ws.Range("BB88").Value = 55
$ws.Range("BC88").Value = 100
$ws.Range("BD88").Value = 350
